# Generate Report for Handoff
# Adds two new localization entries ("569073cc-b091-48bc-a4f7-e62ee521fd94.md"
# and "9d055233-3e97-4a31-9dde-33f35091be4e.md") to the Overview / zh-cn /
# de-de report sheets, each "Ready for handoff".

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

function Set-RowValues($sheet, $rowNum, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $sheet.Cells.Item($rowNum, $i + 1).Value = $values[$i]
    }
}

# ---------------------------------------------------------------------------
# Overview sheet: shift the existing "c3d12459" row down, insert the two new
# files above it (rows 3 and 4), and keep c3d12459 as row 5.
# ---------------------------------------------------------------------------

Set-RowValues $wsOverview 3 @(
    "569073cc-b091-48bc-a4f7-e62ee521fd94.md",
    "e2e\569073cc-b091-48bc-a4f7-e62ee521fd94.md",
    ".md",
    "",
    "Ready for handoff",
    "Ready for handoff",
    "2016-09-02 20:47:35"
)

Set-RowValues $wsOverview 4 @(
    "9d055233-3e97-4a31-9dde-33f35091be4e.md",
    "e2e\9d055233-3e97-4a31-9dde-33f35091be4e.md",
    ".md",
    "",
    "Ready for handoff",
    "Ready for handoff",
    "2016-09-02 20:47:35"
)

Set-RowValues $wsOverview 5 @(
    "c3d12459-73c3-4368-afdd-66ba1d0eb845.md",
    "e2e\c3d12459-73c3-4368-afdd-66ba1d0eb845.md",
    ".md",
    "",
    "Ready for handoff",
    "Ready for handoff",
    "2016-09-02 20:46:22"
)

# date-style column (matches the existing "Latest HO Xliff Generate Date" column format)
$wsOverview.Range("G3:G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/783257e1b99f92e75aa4b82b2c96585368467f5e/e2e/569073cc-b091-48bc-a4f7-e62ee521fd94.md", "", "", "e2e\569073cc-b091-48bc-a4f7-e62ee521fd94.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/86aa02f2abc4a4a5821b9e386c077a58613b108c/e2e/9d055233-3e97-4a31-9dde-33f35091be4e.md", "", "", "e2e\9d055233-3e97-4a31-9dde-33f35091be4e.md") | Out-Null

# The hyperlink that used to live on B3 (c3d12459) now belongs on B5; recreate
# it there (Excel drops a cell's old hyperlink when the cell's row shifts
# content instead of physically moving rows, so it has to be re-added).
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/75404f1344a52ee84cd7526cc9a8dcf722adcf0b/e2e/c3d12459-73c3-4368-afdd-66ba1d0eb845.md", "", "", "e2e\c3d12459-73c3-4368-afdd-66ba1d0eb845.md") | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G5")) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------

Set-RowValues $wsZhCn 3 @(
    "569073cc-b091-48bc-a4f7-e62ee521fd94.md",
    ".md",
    "Ready for handoff",
    "e2e",
    "ht",
    "False",
    "569073cc-b091-48bc-a4f7-e62ee521fd94.783257e1b99f92e75aa4b82b2c96585368467f5e.zh-cn.xlf",
    "2016-09-02 20:47:30",
    "",
    "",
    "0001-01-01 00:00:00",
    "",
    "True",
    "",
    "False",
    ""
)

Set-RowValues $wsZhCn 4 @(
    "9d055233-3e97-4a31-9dde-33f35091be4e.md",
    ".md",
    "Ready for handoff",
    "e2e",
    "ht",
    "False",
    "9d055233-3e97-4a31-9dde-33f35091be4e.86aa02f2abc4a4a5821b9e386c077a58613b108c.zh-cn.xlf",
    "2016-09-02 20:47:30",
    "",
    "",
    "0001-01-01 00:00:00",
    "",
    "True",
    "",
    "False",
    ""
)

Set-RowValues $wsZhCn 5 @(
    "c3d12459-73c3-4368-afdd-66ba1d0eb845.md",
    ".md",
    "Ready for handoff",
    "e2e",
    "ht",
    "False",
    "c3d12459-73c3-4368-afdd-66ba1d0eb845.df08eefcfb5bf83e82e29c7cc6ed7e1246956f74.zh-cn.xlf",
    "2016-09-02 20:46:17",
    "",
    "",
    "0001-01-01 00:00:00",
    "",
    "True",
    "",
    "False",
    ""
)

$wsZhCn.Range("H3:H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K3:K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/783257e1b99f92e75aa4b82b2c96585368467f5e/e2e/569073cc-b091-48bc-a4f7-e62ee521fd94.md", "", "", "569073cc-b091-48bc-a4f7-e62ee521fd94.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/86aa02f2abc4a4a5821b9e386c077a58613b108c/e2e/9d055233-3e97-4a31-9dde-33f35091be4e.md", "", "", "9d055233-3e97-4a31-9dde-33f35091be4e.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/75404f1344a52ee84cd7526cc9a8dcf722adcf0b/e2e/c3d12459-73c3-4368-afdd-66ba1d0eb845.md", "", "", "c3d12459-73c3-4368-afdd-66ba1d0eb845.md") | Out-Null

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P5")) | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------

Set-RowValues $wsDeDe 3 @(
    "569073cc-b091-48bc-a4f7-e62ee521fd94.md",
    ".md",
    "Ready for handoff",
    "e2e",
    "ht",
    "False",
    "569073cc-b091-48bc-a4f7-e62ee521fd94.783257e1b99f92e75aa4b82b2c96585368467f5e.de-de.xlf",
    "2016-09-02 20:47:35",
    "",
    "",
    "0001-01-01 00:00:00",
    "",
    "True",
    "",
    "False",
    ""
)

Set-RowValues $wsDeDe 4 @(
    "9d055233-3e97-4a31-9dde-33f35091be4e.md",
    ".md",
    "Ready for handoff",
    "e2e",
    "ht",
    "False",
    "9d055233-3e97-4a31-9dde-33f35091be4e.86aa02f2abc4a4a5821b9e386c077a58613b108c.de-de.xlf",
    "2016-09-02 20:47:35",
    "",
    "",
    "0001-01-01 00:00:00",
    "",
    "True",
    "",
    "False",
    ""
)

Set-RowValues $wsDeDe 5 @(
    "c3d12459-73c3-4368-afdd-66ba1d0eb845.md",
    ".md",
    "Ready for handoff",
    "e2e",
    "ht",
    "False",
    "c3d12459-73c3-4368-afdd-66ba1d0eb845.df08eefcfb5bf83e82e29c7cc6ed7e1246956f74.de-de.xlf",
    "2016-09-02 20:46:22",
    "",
    "",
    "0001-01-01 00:00:00",
    "",
    "True",
    "",
    "False",
    ""
)

$wsDeDe.Range("H3:H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K3:K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/783257e1b99f92e75aa4b82b2c96585368467f5e/e2e/569073cc-b091-48bc-a4f7-e62ee521fd94.md", "", "", "569073cc-b091-48bc-a4f7-e62ee521fd94.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/86aa02f2abc4a4a5821b9e386c077a58613b108c/e2e/9d055233-3e97-4a31-9dde-33f35091be4e.md", "", "", "9d055233-3e97-4a31-9dde-33f35091be4e.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/75404f1344a52ee84cd7526cc9a8dcf722adcf0b/e2e/c3d12459-73c3-4368-afdd-66ba1d0eb845.md", "", "", "c3d12459-73c3-4368-afdd-66ba1d0eb845.md") | Out-Null

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P5")) | Out-Null
